$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$sub3 = [string][char]0x2083

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '56.278.93'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -2.06%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.373.50'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -1.48%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '501.49'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -1.16%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '129.84'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -2.34%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.51%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.546'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -2.05%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.380.03'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -2.62%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0981'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.18%  '

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.45%  '

$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = 'Toncoin'
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.78'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +3.82%  '

$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = 'Cardano'
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.324'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.55%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.797.63'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -1.58%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '56.264.08'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.80%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.50'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -1.63%  '

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.75%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.367.27'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -3.92%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.02'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -2.66%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.03'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -1.91%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '306.88'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -2.20%  '

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -2.32%  '

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.16%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '65.43'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.41%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.997'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.44%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.367'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -3.50%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.148'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -3.70%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.24'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -4.31%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '171.54'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -1.29%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0' + $sub3 + '0714'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -2.48%  '

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -2.69%  '

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.24%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.75'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -6.74%  '

$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = 'FirstDigitalUSD'
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.998'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.60%  '

$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = 'Fetch.AI'
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.08'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -3.92%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '17.55'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -2.22%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.76'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -1.57%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '36.07'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -1.57%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.795'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -2.24%  '

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -5.59%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '130.59'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -3.32%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.35'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.89%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.73'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -5.38%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.563'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -1.32%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0906'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -1.13%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '240.34'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -5.81%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0480'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -2.50%  '

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -2.67%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '17.06'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.27%  '

$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = 'BitgetToken'
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = 'https://coinranking.com/coin/q7gMmMdLb+bitgettoken-bgb'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.951'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.55%  '
